$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.364.37'
$ws.Range("E2").Value = '  +1.61%  '
$ws.Range("D3").Value = '1.826.08'
$ws.Range("E3").Value = '  +2.86%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'317.34"
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'0.4057"
$ws.Range("E8").Value = '  +8.70%  '
$ws.Range("D9").Value = "'0.07602"
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").Value = "'1.105"
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("D12").Value = "'6.328"
$ws.Range("E12").Value = '  +4.55%  '
$ws.Range("D13").Value = "'1.001"
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("D14").Value = "'7.609"
$ws.Range("E14").Value = '  +5.75%  '
$ws.Range("E15").Value = '  +1.79%  '
$ws.Range("D16").Value = '1.824.20'
$ws.Range("E16").Value = '  +3.17%  '
$ws.Range("D17").Value = "'89.32"
$ws.Range("E17").Value = '  +1.86%  '
$ws.Range("E18").Value = '  +2.19%  '
$ws.Range("D19").Value = "'0.06607"
$ws.Range("E19").Value = '  +3.17%  '
$ws.Range("D20").Value = "'17.63"
$ws.Range("E20").Value = '  +2.33%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = "'6.106"
$ws.Range("E22").Value = '  +3.98%  '
$ws.Range("D23").Value = '28.378.41'
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("D25").Value = "'2.184"
$ws.Range("E25").Value = '  +5.30%  '
$ws.Range("D26").Value = "'2.468"
$ws.Range("D27").Value = "'157.84"
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("D28").Value = "'20.55"
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("D29").Value = '2.035.63'
$ws.Range("E29").Value = '  +3.25%  '
$ws.Range("D30").Value = "'123.94"
$ws.Range("E30").Value = '  +3.66%  '
$ws.Range("E31").Value = '  +1.27%  '
$ws.Range("D32").Value = "'0.1096"
$ws.Range("E32").Value = '  +5.57%  '
$ws.Range("D33").Value = "'5.657"
$ws.Range("E33").Value = '  +3.06%  '
$ws.Range("D34").Value = "'3.645"
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("D35").Value = "'0.07297"
$ws.Range("E35").Value = '  +15.24%  '
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("D37").Value = "'0.02345"
$ws.Range("E37").Value = '  +3.87%  '
$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").Value = "'8.881"
$ws.Range("E38").Value = '  +6.33%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = "'5.189"
$ws.Range("E39").Value = '  +4.88%  '
$ws.Range("D40").Value = "'0.6254"
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("E41").Value = '  +2.98%  '
$ws.Range("D42").Value = "'1.186"
$ws.Range("E42").Value = '  +1.46%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = "'1.399"
$ws.Range("E44").Value = '  -1.88%  '
$ws.Range("E45").Value = '  +2.19%  '
$ws.Range("D46").Value = "'3.706"
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("E47").Value = '  +1.99%  '
$ws.Range("D48").Value = "'125.37"
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = "'1.986"
$ws.Range("E49").Value = '  +3.63%  '
$ws.Range("D50").Value = "'1.203"
$ws.Range("E50").Value = '  +2.02%  '
$ws.Range("D51").Value = "'0.06889"
$ws.Range("E51").Value = '  +1.55%  '
